$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "C2"  = -0.1474887606676356
    "E2"  = 1.265019766896436
    "C3"  = 0.337821977117625
    "E3"  = 1.858657482882586
    "C4"  = -1.890773121057054
    "E4"  = -3.246097549514837
    "C5"  = 1.89156560755015
    "E5"  = 1.609625625599986
    "C6"  = 1.113165545862094
    "E6"  = 1.609625625600009
    "C7"  = 1.070385798714391
    "E7"  = 3.238605209599998
    "C8"  = 1.384186838979806
    "E8"  = 2.777885851461503
    "C9"  = 2.349355943833098
    "E9"  = 2.436566844071941
    "C10" = 1.786425635558397
    "E10" = 1.694971351092267
    "C11" = 1.331333081915509
    "E11" = 1.216098605743365
    "C12" = 1.282262557986447
    "E12" = 1.784618024189011
    "C13" = 2.247109253368285
    "E13" = 4.887093273600018
    "C14" = -4.247034401476801
    "E14" = -12.1986023424
    "C15" = -2.608215948579529
    "E15" = 12.1815000816919
    "C16" = 4.863085601670813
    "E16" = 5.870037016039187
    "C17" = -1.44371442952016
    "E17" = -1.24582517146522
    "C18" = 0.06625622369935691
    "E18" = -0.1040473946152809
    "C19" = 0.9919038146506631
    "E19" = 0.8749339604052775
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
